$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row 1, shifting all existing data down by one row.
$ws.Rows.Item(1).Insert()

# Put the new header/notice text in A1.
$ws.Range("A1").Value = "If you can not find the data which you want,please contact Brad.Qiu"

# Merge A1:E1 into a single banner cell.
$ws.Range("A1:E1").Merge()

# Center + wrap the whole merged banner range.
$ws.Range("A1:E1").HorizontalAlignment = -4108
$ws.Range("A1:E1").WrapText = $true

# Make the banner text stand out: big, red font (only needs to live on A1,
# the merge anchor cell).
$ws.Range("A1").Font.Color = 255
$ws.Range("A1").Font.Size = 18

# Give the banner row extra height so the wrapped text fits.
$ws.Rows.Item(1).RowHeight = 56.25

# Match the recorded selection from the authored workbook.
$ws.Range("F17").Select()
